$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the "MainPump Driver" worksheet to "MainDriver".
#    Excel cascades this rename into every defined name and every
#    plain (non-chart) formula that referenced the old sheet name,
#    e.g. driver_name / gear_ratio / power_curve defined names and the
#    B10 formulas on UWPump / MainPump.
# ---------------------------------------------------------------------
$wsDriver = $wb.Worksheets.Item("MainPump Driver")
$wsDriver.Name = "MainDriver"

$wsUW   = $wb.Worksheets.Item("UWPump")
$wsMain = $wb.Worksheets.Item("MainPump")

# ---------------------------------------------------------------------
# 2. The rename-cascade does not always requote/clean up every single
#    formula reference, so make the post-rename formulas explicit to
#    match the target workbook exactly.
# ---------------------------------------------------------------------
$wsUW.Range("B10").Formula   = "=MainDriver!B14"
$wsMain.Range("B10").Formula = "=MainDriver!B14"

# ---------------------------------------------------------------------
# 3. UWPump!J4 used to show the driver name ('MainPump Driver'!driver_name);
#    it now shows the pump name instead.
# ---------------------------------------------------------------------
$wsUW.Range("J4").Formula = "=MainPump!pump_name"

# ---------------------------------------------------------------------
# 3b. The scatter chart that lives on the MainDriver sheet itself also
#     referenced the old sheet name in its series formula; repoint the
#     series name reference at the renamed sheet as well.
# ---------------------------------------------------------------------
try {
    $wsDriverForChart = $wb.Worksheets.Item("MainDriver")
    $chartObj = $wsDriverForChart.ChartObjects(1)
    $series1 = $chartObj.Chart.SeriesCollection(1)
    $series1.Formula = "=SERIES(MainDriver!`$B`$5,MainDriver!`$A`$6:`$A`$14,MainDriver!`$B`$6:`$B`$14,1)"
} catch {
    Write-Host "chart series formula update skipped: $_"
}

# ---------------------------------------------------------------------
# 4. Selection / active-tab bookkeeping:
#    - UWPump loses the "tabSelected" flag and its cursor moves to B2.
#    - MainDriver (the renamed sheet) becomes the active/selected tab,
#      with its cursor left on H9 (unchanged).
# ---------------------------------------------------------------------
$wsUW.Range("B2").Select() | Out-Null

$wsDriverAgain = $wb.Worksheets.Item("MainDriver")
$wsDriverAgain.Activate() | Out-Null
